$d = $word.ActiveDocument

# 1. Student name: "Aliza Hanum Anggani" -> "Reo Sahobby"
$d.Content.Find.Execute("Aliza Hanum Anggani", $false, $false, $false, $false, $false, $true, 1, $false, "Reo Sahobby", 2)

# 2. NIM score digits: "123170090" -> "123170067"
$d.Content.Find.Execute("123170090", $false, $false, $false, $false, $false, $true, 1, $false, "123170067", 2)

# 3. Judul Tugas Akhir - first line
$d.Content.Find.Execute("Penerapan Metode Ensemble Untuk Mengatasi ", $false, $false, $false, $false, $false, $true, 1, $false, "Prediksi Penyakit Jantung dengan Menggunakan ", 2)

# 4. Judul Tugas Akhir - second line
$d.Content.Find.Execute("Ketidakseimbangan Dataset Pada Klasifikasi Penyakit ", $false, $false, $false, $false, $false, $true, 1, $false, "Algoritma XgBoost dan Randomized Search Optimizer ", 2)

# 5. Judul Tugas Akhir - third line removed entirely
$d.Content.Find.Execute("Kesehatan Mental Menggunakan Algoritma Naïve Bayes", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
